# "11 May - Noche"
# Corrects the Materia (and, where applicable, Calificacion) values that had
# been swapped between pairs of rows belonging to the same student on
# several "Rescatables" sheets.

$wb = $excel.ActiveWorkbook

function Swap-ERows {
    param(
        [string]$SheetName,
        [int]$Row1,
        [int]$Row2,
        [bool]$SwapF
    )

    $ws = $wb.Worksheets.Item($SheetName)

    $e1 = $ws.Cells.Item($Row1, 5).Value()
    $e2 = $ws.Cells.Item($Row2, 5).Value()
    $ws.Cells.Item($Row1, 5).Value = $e2
    $ws.Cells.Item($Row2, 5).Value = $e1

    if ($SwapF) {
        $f1 = $ws.Cells.Item($Row1, 6).Value()
        $f2 = $ws.Cells.Item($Row2, 6).Value()
        $ws.Cells.Item($Row1, 6).Value = $f2
        $ws.Cells.Item($Row2, 6).Value = $f1
    }
}

# 4AEM: rows 18/19 swap Materia + Calificacion; rows 21/22 swap Materia only
Swap-ERows "4AEM" 18 19 $true
Swap-ERows "4AEM" 21 22 $false

# 4BEM: rows 2/3 and rows 8/9 swap Materia only
Swap-ERows "4BEM" 2 3 $false
Swap-ERows "4BEM" 8 9 $false

# 6ASM: rows 4/5 swap Materia + Calificacion
Swap-ERows "6ASM" 4 5 $true

# 4AEV: rows 2/3 swap Materia + Calificacion
Swap-ERows "4AEV" 2 3 $true

# 6AEV: rows 4/5 swap Materia + Calificacion
Swap-ERows "6AEV" 4 5 $true
